# Remove the "Solution" column (header-only column D) from the sheet.
# The header row (row 1) had 7 columns: Question No., Problem Name,
# Problem Statement, Solution, Technique, Topic, Difficulty. The data row
# (row 2) only ever had 6 values (it never had a "Solution" value), so
# removing the "Solution" header shifts the remaining row-1 headers left
# and drops the now-empty trailing column, while row 2's data is left
# exactly as it was.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Technique"
$ws.Range("E1").Value = "Topic"
$ws.Range("F1").Value = "Difficulty"
$ws.Range("G1").ClearContents()
